# Rename the worksheet tab (Geo 121 -> BUS 333) and refresh the three
# "Course Color" RGB component cells (D1/E1/F1) with their new values.
#
# D1/E1/F1 already hold text (shared-string) cells for the color
# components, not numbers - simply assigning Range.Value with a numeric-
# looking string would silently re-type the cell as a Number. To keep
# them as text (matching the existing column formatting/behavior) we
# build the literal through a formula and then freeze it back down to a
# plain value with Copy + PasteSpecial(xlPasteValues), exactly like using
# Paste Values after typing ="204" in the formula bar.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "BUS 333"

function Set-TextValue($range, [string]$text) {
    $cell = $ws.Range($range)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

Set-TextValue "E1" "255"
Set-TextValue "F1" "255"
Set-TextValue "D1" "204"

$excel.CutCopyMode = $false
